# Auto-generated script to apply numeric corrections to Pandaemonium_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 303752.38
$ws.Range("I64").Value = 489222.34
$ws.Range("K64").Value = 489222.34
$ws.Range("M64").Value = -488974.34
$ws.Range("H67").Value = 303752.38
$ws.Range("I67").Value = 489222.34
$ws.Range("K67").Value = 489222.34
$ws.Range("M67").Value = -488364.34
$ws.Range("H80").Value = 7004.7646
$ws.Range("I80").Value = 531.2
$ws.Range("J80").Value = 16252.714
$ws.Range("K80").Value = 1593.6
$ws.Range("L80").Value = 48758.142
$ws.Range("M80").Value = -595.6000000000001
$ws.Range("N80").Value = -50754.142
$ws.Range("H83").Value = 7004.7646
$ws.Range("I83").Value = 531.2
$ws.Range("J83").Value = 16252.714
$ws.Range("K83").Value = 4780.8
$ws.Range("L83").Value = 146274.426
$ws.Range("M83").Value = 211.1999999999998
$ws.Range("N83").Value = -156258.426
$ws.Range("H100").Value = 1427.8125
$ws.Range("I100").Value = 1435.3
$ws.Range("J100").Value = 1415.3334
$ws.Range("K100").Value = 1435.3
$ws.Range("L100").Value = 1415.3334
$ws.Range("M100").Value = -894.3
$ws.Range("N100").Value = -2497.3334
$ws.Range("H129").Value = 847.1111
$ws.Range("I129").Value = 286.75
$ws.Range("J129").Value = 944.56525
$ws.Range("K129").Value = 860.25
$ws.Range("L129").Value = 2833.69575
$ws.Range("M129").Value = 4139.75
$ws.Range("N129").Value = -12833.69575
$ws.Range("H132").Value = 1935.5294
$ws.Range("I132").Value = 1994.0312
$ws.Range("J132").Value = 999.5
$ws.Range("K132").Value = 5982.0936
$ws.Range("L132").Value = 2998.5
$ws.Range("M132").Value = -3452.0936
$ws.Range("N132").Value = -8058.5
$ws.Range("H137").Value = 2350.0222
$ws.Range("I137").Value = 1646.5862
$ws.Range("J137").Value = 3625
$ws.Range("K137").Value = 4939.7586
$ws.Range("L137").Value = 10875
$ws.Range("M137").Value = -2389.7586
$ws.Range("N137").Value = -15975
$ws.Range("H138").Value = 3364.648
$ws.Range("I138").Value = 1767.4073
$ws.Range("J138").Value = 4344.773
$ws.Range("K138").Value = 5302.2219
$ws.Range("L138").Value = 13034.319
$ws.Range("M138").Value = -162.2219000000005
$ws.Range("N138").Value = -23314.319

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16371.56
$ws.Range("I32").Value = 16909.973
$ws.Range("J32").Value = 11884.777
$ws.Range("K32").Value = 16909.973
$ws.Range("L32").Value = 11884.777
$ws.Range("M32").Value = -16622.973
$ws.Range("N32").Value = -12458.777
$ws.Range("H132").Value = 4540.4653
$ws.Range("I132").Value = 1688.2354
$ws.Range("J132").Value = 15315.556
$ws.Range("K132").Value = 5064.706200000001
$ws.Range("L132").Value = 45946.66800000001
$ws.Range("M132").Value = -2534.706200000001
$ws.Range("N132").Value = -51006.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 34601
$ws.Range("J110").Value = 34601
$ws.Range("L110").Value = 34601
$ws.Range("N110").Value = -42781
$ws.Range("H112").Value = 26232.736
$ws.Range("J112").Value = 26232.736
$ws.Range("L112").Value = 26232.736
$ws.Range("N112").Value = -29186.736
$ws.Range("H114").Value = 79800
$ws.Range("J114").Value = 79800
$ws.Range("L114").Value = 79800
$ws.Range("N114").Value = -88478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4589.075
$ws.Range("I31").Value = 4504.839
$ws.Range("J31").Value = 4879.222
$ws.Range("K31").Value = 4504.839
$ws.Range("L31").Value = 4879.222
$ws.Range("M31").Value = -4209.839
$ws.Range("N31").Value = -5469.222
$ws.Range("H34").Value = 4589.075
$ws.Range("I34").Value = 4504.839
$ws.Range("J34").Value = 4879.222
$ws.Range("K34").Value = 4504.839
$ws.Range("L34").Value = 4879.222
$ws.Range("M34").Value = -4302.839
$ws.Range("N34").Value = -5283.222
$ws.Range("H132").Value = 4298.1665
$ws.Range("I132").Value = 5189.2856
$ws.Range("J132").Value = 3050.6
$ws.Range("K132").Value = 15567.8568
$ws.Range("L132").Value = 9151.799999999999
$ws.Range("M132").Value = -13037.8568
$ws.Range("N132").Value = -14211.8
$ws.Range("H134").Value = 3065.12
$ws.Range("I134").Value = 2022.8889
$ws.Range("J134").Value = 4288.609
$ws.Range("K134").Value = 6068.6667
$ws.Range("L134").Value = 12865.827
$ws.Range("M134").Value = -3533.6667
$ws.Range("N134").Value = -17935.827

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 13060
$ws.Range("J97").Value = 22478
$ws.Range("L97").Value = 67434
$ws.Range("N97").Value = -68426
$ws.Range("H115").Value = 2350
$ws.Range("I115").Value = 1400
$ws.Range("J115").Value = 2666.6667
$ws.Range("K115").Value = 4200
$ws.Range("L115").Value = 8000.000100000001
$ws.Range("M115").Value = -3025
$ws.Range("N115").Value = -10350.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 17930
$ws.Range("J136").Value = 17930
$ws.Range("L136").Value = 53790
$ws.Range("N136").Value = -58890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6472.1377
$ws.Range("I132").Value = 7014.7
$ws.Range("J132").Value = 5266.4443
$ws.Range("K132").Value = 21044.1
$ws.Range("L132").Value = 15799.3329
$ws.Range("M132").Value = -18514.1
$ws.Range("N132").Value = -20859.3329
$ws.Range("H136").Value = 4898.089
$ws.Range("I136").Value = 2973.1333
$ws.Range("J136").Value = 8748
$ws.Range("K136").Value = 8919.3999
$ws.Range("L136").Value = 26244
$ws.Range("M136").Value = -6369.3999
$ws.Range("N136").Value = -31344

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 43429
$ws.Range("J46").Value = 43429
$ws.Range("L46").Value = 43429
$ws.Range("N46").Value = -43891
$ws.Range("H96").Value = 1951
$ws.Range("J96").Value = 1934.6666
$ws.Range("L96").Value = 1934.6666
$ws.Range("N96").Value = -4680.6666
$ws.Range("H132").Value = 2638.6155
$ws.Range("I132").Value = 1128
$ws.Range("J132").Value = 3746.4
$ws.Range("K132").Value = 3384
$ws.Range("L132").Value = 11239.2
$ws.Range("M132").Value = -854
$ws.Range("N132").Value = -16299.2
$ws.Range("H134").Value = 43429
$ws.Range("J134").Value = 43429
$ws.Range("L134").Value = 130287
$ws.Range("N134").Value = -135357
$ws.Range("H136").Value = 6252.091
$ws.Range("I136").Value = 5330.2
$ws.Range("J136").Value = 7358.36
$ws.Range("K136").Value = 15990.6
$ws.Range("L136").Value = 22075.08
$ws.Range("M136").Value = -13440.6
$ws.Range("N136").Value = -27175.08
